$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# Row 2
Set-TextValue "D2" "67.685.70"
$ws.Range("E2").Value = "  -1.80%  "

# Row 3
Set-TextValue "D3" "3.266.52"
$ws.Range("E3").Value = "  -0.95%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
Set-TextValue "D5" "580.31"
$ws.Range("E5").Value = "  -1.06%  "

# Row 6
Set-TextValue "D6" "184.61"
$ws.Range("E6").Value = "  +0.58%  "

# Row 8
Set-TextValue "D8" "0.602"
$ws.Range("E8").Value = "  +0.54%  "

# Row 9
Set-TextValue "D9" "0.130"
$ws.Range("E9").Value = "  -4.41%  "

# Row 10
$ws.Range("E10").Value = "  -1.37%  "

# Row 11
Set-TextValue "D11" "0.408"
$ws.Range("E11").Value = "  -3.50%  "

# Row 12
Set-TextValue "D12" "3.821.24"
$ws.Range("E12").Value = "  -1.40%  "

# Row 13
$ws.Range("E13").Value = "  +0.50%  "

# Row 14
Set-TextValue "D14" "27.40"
$ws.Range("E14").Value = "  -5.88%  "

# Row 15
Set-TextValue "D15" "67.750.22"
$ws.Range("E15").Value = "  -1.67%  "

# Row 16
Set-TextValue "D16" "0.0000168"
$ws.Range("E16").Value = "  -2.49%  "

# Row 17
Set-TextValue "D17" "3.235.03"
$ws.Range("E17").Value = "  -1.22%  "

# Row 18
Set-TextValue "D18" "5.72"
$ws.Range("E18").Value = "  -2.27%  "

# Row 19
Set-TextValue "D19" "13.43"
$ws.Range("E19").Value = "  -1.65%  "

# Row 20
Set-TextValue "D20" "396.91"
$ws.Range("E20").Value = "  +0.66%  "

# Row 21
Set-TextValue "D21" "7.56"
$ws.Range("E21").Value = "  -2.62%  "

# Row 22
$ws.Range("E22").Value = "  +0.19%  "

# Row 23
$ws.Range("E23").Value = "  -1.49%  "

# Row 24
$ws.Range("E24").Value = "  -1.82%  "

# Row 25
$ws.Range("E25").Value = "  -3.85%  "

# Row 26
Set-TextValue "D26" "0.188"
$ws.Range("E26").Value = "  -0.72%  "

# Row 27
Set-TextValue "D27" "9.51"
$ws.Range("E27").Value = "  -2.23%  "

# Row 28
$ws.Range("E28").Value = "  +0.26%  "

# Row 29
$ws.Range("E29").Value = "  -2.21%  "

# Row 30
Set-TextValue "D30" "22.63"
$ws.Range("E30").Value = "  -2.24%  "

# Row 31
Set-TextValue "D31" "5.46"
$ws.Range("E31").Value = "  -5.40%  "

# Row 32
Set-TextValue "D32" "6.93"
$ws.Range("E32").Value = "  -3.55%  "

# Row 34
$ws.Range("E34").Value = "  -4.87%  "

# Row 35
Set-TextValue "D35" "164.72"
$ws.Range("E35").Value = "  +0.70%  "

# Row 36
Set-TextValue "D36" "1.46"
$ws.Range("E36").Value = "  -5.10%  "

# Row 37
Set-TextValue "D37" "1.89"
$ws.Range("E37").Value = "  -0.64%  "

# Row 38
Set-TextValue "D38" "26.86"
$ws.Range("E38").Value = "  +1.60%  "

# Row 39
Set-TextValue "D39" "0.807"
$ws.Range("E39").Value = "  -3.69%  "

# Row 40
Set-TextValue "D40" "4.50"
$ws.Range("E40").Value = "  -2.51%  "

# Row 41
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D41" "2.667.50"
$ws.Range("E41").Value = "  +1.73%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D42" "6.26"
$ws.Range("E42").Value = "  -5.15%  "

# Row 43
Set-TextValue "D43" "40.76"
$ws.Range("E43").Value = "  -2.00%  "

# Row 44
Set-TextValue "D44" "0.0680"
$ws.Range("E44").Value = "  -1.91%  "

# Row 45
Set-TextValue "D45" "2.44"
$ws.Range("E45").Value = "  -6.07%  "

# Row 46
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D46" "335.44"
$ws.Range("E46").Value = "  -2.25%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D47" "24.59"
$ws.Range("E47").Value = "  -2.50%  "

# Row 48
$ws.Range("E48").Value = "  -3.70%  "

# Row 49
Set-TextValue "D49" "6.29"
$ws.Range("E49").Value = "  -0.51%  "

# Row 50
Set-TextValue "D50" "0.101"
$ws.Range("E50").Value = "  -1.71%  "

# Row 51
$ws.Range("E51").Value = "  -2.18%  "
